$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: "id" / "x" / "y"  ->  "index" / "x" / "y"
$ws.Range("A1").Value = "index"
$ws.Range("B1").Value = "x"
$ws.Range("C1").Value = "y"

# Column A used to hold text labels ("id0".."id2"); it now holds the
# numeric row index (0-based) instead.
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2

# Leave the cursor where the author left it when the file was saved.
$ws.Range("C10").Select()
